$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed cryptos list values (price / 1h volume %, and the two
# swapped coin rows) from the latest GitHub Actions data pull.
#
# Price cells in column D that look like plain numbers are entered with a
# leading apostrophe, exactly as a user typing into Excel would do to force
# text, so the COM layer does not silently convert them into numeric cells.
# The style is reset to "Normal" right after so no quote-prefix styling
# lingers on the cell.

$ws.Range("D2").Value = '70.804.64'
$ws.Range("E2").Value = '  +1.97%  '
$ws.Range("D3").Value = '3.540.14'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''607.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.42%  '
$ws.Range("D6").Value = '''172.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = '''0.619'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.11%  '
$ws.Range("D8").Value = '3.537.53'
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("D11").Value = '''6.83'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").Value = '''46.87'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").Value = '4.120.82'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '''8.39'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.76%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '''618.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.51%  '
$ws.Range("D18").Value = '70.728.06'
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("D19").Value = '3.530.76'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '''0.120'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.54%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").Value = '''9.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -16.48%  '
$ws.Range("E24").Value = '  -2.63%  '
$ws.Range("D25").Value = '''96.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.24%  '
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  -1.91%  '
$ws.Range("E29").Value = '  +1.68%  '
$ws.Range("D30").Value = '''9.07'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.91%  '
$ws.Range("D31").Value = '''8.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.90%  '
$ws.Range("D32").Value = '''3.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.27%  '
$ws.Range("B33").Value = 'Mantle'
$ws.Range("C33").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D33").Value = '''1.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '''6.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.28%  '
$ws.Range("D35").Value = '''572.84'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.63%  '
$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").Value = '''3.66'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.42%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '''0.102'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("D38").Value = '''10.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("D39").Value = '''57.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("D40").Value = '''0.0473'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.06%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").Value = '''0.143'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.52%  '
$ws.Range("D43").Value = '3.342.92'
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("D44").Value = '''0.324'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("D45").Value = '''3.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.46%  '
$ws.Range("D46").Value = '''33.07'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("E48").Value = '  +1.51%  '
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("D50").Value = '''133.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").Value = '''5.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.31%  '
